$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.44870138168335
$ws.Range("B1").Value = 5.373844146728516
$ws.Range("C1").Value = 8.447382926940918
$ws.Range("D1").Value = 8.41199779510498
$ws.Range("E1").Value = 3.61477518081665
